$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 currently holds the "Dutch Eredivisie" match, which is being dropped.
# Row 3 holds the "Mexican Liga MX" match, which becomes the new row 2.
# Copy row 3 (text + numbers) over row 2 so text cells (League/Date/Time/Home/Away)
# keep their original inline-string representation instead of being re-parsed
# (e.g. avoiding "2026-01-09" turning into a date serial number).
$ws.Range("A3:AO3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# Now apply the updated odds for this match into row 2.
$ws.Range("F2").Value = 3.25
$ws.Range("G2").Value = 3.35
$ws.Range("H2").Value = 2.52
$ws.Range("I2").Value = 2.6
$ws.Range("J2").Value = 3.25
$ws.Range("K2").Value = 3.35
$ws.Range("L2").Value = 1.61
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 3.1
$ws.Range("O2").Value = 1.46
$ws.Range("P2").Value = 1.72
$ws.Range("Q2").Value = 2.34
$ws.Range("R2").Value = 1.26
$ws.Range("S2").Value = 4.5
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.59
$ws.Range("W2").Value = 1.45
$ws.Range("X2").Value = 10
$ws.Range("Y2").Value = 9.800000000000001
$ws.Range("Z2").Value = 15.5
$ws.Range("AA2").Value = 38
$ws.Range("AB2").Value = 10.5
$ws.Range("AC2").Value = 7.2
$ws.Range("AD2").Value = 13
$ws.Range("AE2").Value = 34
$ws.Range("AF2").Value = 23
$ws.Range("AG2").Value = 15
$ws.Range("AH2").Value = 22
$ws.Range("AI2").Value = 60
$ws.Range("AJ2").Value = 75
$ws.Range("AK2").Value = 55
$ws.Range("AL2").Value = 65
$ws.Range("AM2").Value = 180
$ws.Range("AN2").Value = 55
$ws.Range("AO2").Value = 34

# Remove the now-duplicated row 3.
$ws.Rows("3").Delete()
